$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.136510610580444
$ws.Range("B1").Value = 2.56273078918457
$ws.Range("C1").Value = 2.554456472396851
$ws.Range("D1").Value = 2.845099687576294
$ws.Range("E1").Value = 0.535198450088501
